# Updated cryptos list on Sun Oct  8 08:32:08 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '27.930.89'
$ws.Range("E2").Value = '  +0.04%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.629.11'
$ws.Range("E3").Value = '  -0.70%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.03%  '

# Row 5 - BNB
$ws.Range("D5").Value = "'211.56"
$ws.Range("D5").Style = "Normal"

# Row 6 - XRP
$ws.Range("E6").Value = '  +0.01%  '

# Row 8 - Solana
$ws.Range("D8").Value = "'23.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.94%  '

# Row 9 - Cardano
$ws.Range("E9").Value = '  -1.91%  '

# Row 10 - Dogecoin
$ws.Range("E10").Value = '  -0.20%  '

# Row 11 - TRON
$ws.Range("E11").Value = '  +0.43%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.859.81'
$ws.Range("E12").Value = '  -0.68%  '

# Row 13 - WrappedEther
$ws.Range("D13").Value = '1.628.20'
$ws.Range("E13").Value = '  -0.77%  '

# Row 14 - Polkadot
$ws.Range("E14").Value = '  -1.53%  '

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.83%  '

# Row 16 - Litecoin
$ws.Range("D16").Value = "'65.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.55%  '

# Row 17 - WrappedBTC
$ws.Range("D17").Value = '27.912.29'
$ws.Range("E17").Value = '  +0.01%  '

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'231.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.12%  '

# Row 19 - Chainlink
$ws.Range("E19").Value = '  +1.12%  '

# Row 20 - ShibaInu
$ws.Range("E20").Value = '  +0.02%  '

# Row 21 - Dai
$ws.Range("E21").Value = '  -0.05%  '

# Row 22 - Uniswap
$ws.Range("E22").Value = '  -0.35%  '

# Row 23 - Avalanche
$ws.Range("D23").Value = "'10.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.49%  '

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.48%  '

# Row 25 - Monero
$ws.Range("D25").Value = "'154.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.20%  '

# Row 26 - Cosmos
$ws.Range("E26").Value = '  -0.26%  '

# Row 27 - Stellar
$ws.Range("E27").Value = '  +0.15%  '

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'15.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.70%  '

# Row 29 - BinanceUSD
$ws.Range("E29").Value = '  -0.06%  '

# Row 30 - PancakeSwap
$ws.Range("E30").Value = '  -0.76%  '

# Row 31 - Hedera
$ws.Range("D31").Value = "'0.0483"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.27%  '

# Row 32 - Filecoin
$ws.Range("E32").Value = '  +2.54%  '

# Row 34 - Maker
$ws.Range("D34").Value = '1.402.11'
$ws.Range("E34").Value = '  +0.36%  '

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = '  +0.23%  '

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.76%  '

# Row 37 - HuobiToken
$ws.Range("E37").Value = '  -0.40%  '

# Row 38 - VeChain
$ws.Range("E38").Value = '  +2.04%  '

# Row 39 - ImmutableX
$ws.Range("D39").Value = "'0.557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.20%  '

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "'0.865"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.62%  '

# Row 41 - WEMIXToken
$ws.Range("E41").Value = '  -0.46%  '

# Row 42 - PaxDollar
$ws.Range("E42").Value = '  -0.09%  '

# Row 43 - RenderToken
$ws.Range("D43").Value = "'1.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.00%  '

# Row 44 - was Aave, now FraxShare
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = "'5.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.36%  '

# Row 45 - was FraxShare, now Aave
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'66.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.05%  '

# Row 46 - MXToken
$ws.Range("E46").Value = '  -0.59%  '

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = '1.768.70'

# Row 48 - Quant
$ws.Range("D48").Value = "'88.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.19%  '

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = '  -1.19%  '

# Row 50 - Algorand
$ws.Range("E50").Value = '  +1.14%  '

# Row 51 - Cronos
$ws.Range("E51").Value = '  -0.35%  '
